# The "쿼리2" worksheet is populated by a Power Query that hits a live web
# API for each streamer's monthly cumulative "balloon" total and stamps the
# result with DateTime.LocalNow() (column D, "새로고침시간"). This commit is
# the result of the query being refreshed: every row's total (column C,
# "월별 누적별풍선") and refresh timestamp (column D) move to their new
# values, and the user's active selection ends up on F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new "월별 누적별풍선" (column C) value after the refresh.
# Row 11 (152599) is unchanged by the refresh, so it is left alone.
$newTotals = @(
    @{ Row = 2;  Value = 731973 },
    @{ Row = 3;  Value = 651872 },
    @{ Row = 4;  Value = 533126 },
    @{ Row = 5;  Value = 519289 },
    @{ Row = 6;  Value = 428439 },
    @{ Row = 7;  Value = 383363 },
    @{ Row = 8;  Value = 310183 },
    @{ Row = 9;  Value = 278998 },
    @{ Row = 10; Value = 200078 },
    @{ Row = 12; Value = 80089 }
)

foreach ($entry in $newTotals) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Value
}

# Every data row (2-12) picks up the same new refresh timestamp.
$newRefreshTime = 46015.979279421299
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 4).Value = $newRefreshTime
}

# The user's selection after the refresh completed.
[void]$ws.Range("F4").Select()

Write-Host "Query refresh applied: updated 월별 누적별풍선 / 새로고침시간 for rows 2-12."
